$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate column N ("ELET3": header in row 2 + its values in rows 4-37)
# into new column O, preserving formatting, then relabel the new column's
# header to the newly-added ticker "AXIA6". Row 3 is skipped on purpose --
# column N has no cell there, so column O should not get one either.
$ws.Range("N2").Copy($ws.Range("O2"))
$ws.Range("N4:N37").Copy($ws.Range("O4:O37"))
$ws.Range("O2").Value = "AXIA6"

# Leave the selection on the newly added data range, matching the saved view.
$ws.Range("O4:O37").Select() | Out-Null
